$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the empty "Observaciones" value for the existing product (row 2)
$ws.Range("E2").ClearContents()

# Add new manual product row 3
$ws.Range("A3").Value = "TERM50A"
$ws.Range("B3").Value = "SICA"
$ws.Range("C3").Value = "TERMICA 50a"
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = "ferreteria_general"

# Add new manual product row 4
$ws.Range("A4").Value = "TERM32A"
$ws.Range("B4").Value = "JELUZ"
$ws.Range("C4").Value = "TERMICA 32A JELUZ"
$ws.Range("D4").Value = 5000
$ws.Range("E4").Value = "Producto agregado manualmente"
$ws.Range("F4").Value = "ferreteria_general"
